# Regenerate s_vals data to filter save games.
# Updates the numeric values in columns B:G for rows 2-6 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1554434735375247, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 2.461472421528573)
    3 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126)
    4 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    5 = @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 1.896700893398075)
    6 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column B is index 2, C is 3, ... G is 7
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
